$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name
$ws.Range("C3").Value = "Sukhtab Singh Warya"

# Fill in test case data: Condition being Tested (E), Method Inputs (F), Expected Result (G)
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = "59002635, 5550, 1200.00, date.today(), 2.00"
$ws.Range("G7").Value = "Account number: 59002635, Client number: 5550, Balance: 1200.00, Date created: today, Management fee: 2.00"

$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = '59002635, 5550, 1200.00, date.today(), "invalid"'
$ws.Range("G8").Value = "Management fee is set to 2.55"

$ws.Range("E9").Value = "Account created > 10 years ago"
$ws.Range("F9").Value = "59002635, 5550, 1200.00, date.today() - timedelta(days=10*365.26)"
$ws.Range("G9").Value = "Service charge equals BASE_SERVICE_CHARGE"

$ws.Range("E10").Value = "Account created exactly 10 years ago"
$ws.Range("F10").Value = "59002635, 5550, 1200.00, date.today() - timedelta(days=10*365.25)"
$ws.Range("G10").Value = "Service charge equals BASE_SERVICE_CHARGE"

$ws.Range("E11").Value = "Account created < 10 years ago"
$ws.Range("F11").Value = "59002635, 5550, 1200.00, date.today(), 2.00"
$ws.Range("G11").Value = "Service charge equals BASE_SERVICE_CHARGE + 2.00"

$ws.Range("E12").Value = "Account created > 10 years ago"
$ws.Range("F12").Value = "59002635, 5550, 1200.00, date.today() - timedelta(days=10*365.26)"
$ws.Range("G12").Value = "String includes ""Management Fee: Waived"""

$ws.Range("E13").Value = "Account created < 10 years ago"
$ws.Range("F13").Value = "59002635, 5550, 1200.00, date.today(), 2.00"
$ws.Range("G13").Value = "String includes ""Management Fee: $2.00"""

# The now-completed rows (7-13) lose their bold "placeholder" emphasis, matching column A
# being swept into the formatted range as well.
for ($r = 7; $r -le 13; $r++) {
    for ($c = 1; $c -le 7; $c++) {
        $ws.Cells.Item($r, $c).Font.Bold = $false
    }
}
